# "the cell replace engine for dungeon story"
# Scene.xlsx: row 37 (scene 13020002) incorrectly carried a leftover
# "bossunicorn;30" quest-random value in column H (QuestRandom). Clear it
# so the cell goes back to blank, matching the other rows around it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scene")

# Clear the stray value out of H37 (QuestRandom column for scene 13020002).
$ws.Range("H37").ClearContents()

# Leave the cursor on the cell that was just fixed, scrolled so row 13 is
# the first visible row (matches the author's viewport when they made
# the fix).
$ws.Range("H37").Select()
$wn = $excel.ActiveWindow
$wn.ScrollRow = 13
$wn.ScrollColumn = 1
